$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name/title (date moved from 03-28 to 03-29)
$ws.Name = "Through 2022-03-29"

# Update the "March" label in column A row 4
$ws.Range("A4").Value = "March (through 03-29)"

# Update March row (row 4) values
$ws.Range("B4").Value = 28
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 54
$ws.Range("E4").Value = 58
$ws.Range("F4").Value = 29
$ws.Range("G4").Value = 54
$ws.Range("I4").Value = 125

# Update Total row (row 5) values
$ws.Range("B5").Value = 65
$ws.Range("C5").Value = 127
$ws.Range("D5").Value = 185
$ws.Range("E5").Value = 195
$ws.Range("F5").Value = 108
$ws.Range("G5").Value = 195
$ws.Range("I5").Value = 425
